# Update view-count figures (column F) on the "展览" and "全部类型" sheets.
# These numbers were refreshed by the site's data scraper:
#   14052 -> 14055
#   1047  -> 1048
#   14942 -> 14943
#   5442  -> 5444
#   56    -> 57

$wb = $excel.ActiveWorkbook

$updates = @(
    @{ Sheet = "展览";     Cell = "F2";  Value = 14055 },
    @{ Sheet = "展览";     Cell = "F7";  Value = 1048 },
    @{ Sheet = "展览";     Cell = "F9";  Value = 14943 },
    @{ Sheet = "展览";     Cell = "F28"; Value = 5444 },
    @{ Sheet = "展览";     Cell = "F29"; Value = 57 },
    @{ Sheet = "全部类型"; Cell = "F2";  Value = 14055 },
    @{ Sheet = "全部类型"; Cell = "F8";  Value = 1048 },
    @{ Sheet = "全部类型"; Cell = "F10"; Value = 14943 },
    @{ Sheet = "全部类型"; Cell = "F30"; Value = 5444 },
    @{ Sheet = "全部类型"; Cell = "F31"; Value = 57 }
)

foreach ($u in $updates) {
    $ws = $wb.Worksheets.Item($u.Sheet)
    $ws.Range($u.Cell).Value = $u.Value
}
